$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: always search with MatchCase=$true to dodge the Turkish dotted
# capital I (U+0130) case-folding bug in the Find engine (lower-casing it
# expands to two chars and throws off match-range length).
# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# ---------------------------------------------------------------------------

# 1) "Dersin İşlenişi: Mblock ..." -> insert "Öğretmen " right before "Mblock"
$rng = $d.Content
$found = $rng.Find.Execute("Dersin İşlenişi: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.InsertAfter("Öğretmen ")
}

# 2) " programı açılır." -> " programını açar."
$rng = $d.Content
$found = $rng.Find.Execute(" programı açılır.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = " programını açar."
}

# 3) "Arka plan ve karakter eklenir." -> "Arka plan ve karakteri ekler."
$rng = $d.Content
$found = $rng.Find.Execute("Arka plan ve karakter eklenir.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Arka plan ve karakteri ekler."
}

# 4) "Uzantılar kısmından “makine öğrenimi uzantısı” alınır." keeps the same
#    visible text in the diff (only split into extra runs), so nothing to do.

# 5) Remove the whole "Eğitim modeli butonuna basılarak Ok görevi görecek bir
#    araç sağ-sol-yukarı-aşağı şeklinde programa tanıtılır. " paragraph text.
$rng = $d.Content
$found = $rng.Find.Execute("Eğitim modeli butonuna basılarak Ok görevi görecek bir araç sağ-sol-yukarı-aşağı şeklinde programa tanıtılır. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = ""
}

# 6) "Engeller tek bir karakter olarak eklenir." gets new text prepended and a
#    "_GoBack" bookmark inserted right before the final period.
$rng = $d.Content
$found = $rng.Find.Execute("Engeller tek bir karakter olarak eklenir.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Engeller tek bir karakter olarak eklenir."
    $periodStart = $rng.End - 1
    $bmRng = $d.Range($periodStart, $periodStart)
    $d.Bookmarks.Add("_GoBack", $bmRng)
    $prefixRng = $d.Range($rng.Start, $rng.Start)
    $prefixRng.InsertBefore("Eğitim modeli butonuna basılarak sağ-sol-yukarı-aşağı yazılı kağıtlar ve boş ekran programa tanıtır. ")
}

# 7) "Gerekli kodlar yazılır." -> "Öğretmen gerekli kodları yazar." and the
#    old "_GoBack" bookmark located here is removed (it moved to step 6).
$rng = $d.Content
$found = $rng.Find.Execute("Gerekli kodlar yazılır.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "Öğretmen gerekli kodları yazar."
}
